$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Clear the "Pass"/"Fail" DSL markers that are no longer used in column J
$ws.Range("J2").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("J4").ClearContents()

# Move the active selection to E13, matching the edited file's view state
$ws.Activate()
$ws.Range("E13").Select()

# Column J's best-fit width shrinks now that it no longer holds "Pass"/"Fail"
$ws.Columns("J").ColumnWidth = 5.14
